# lab2.docx edit: add "(note that it is on a custom personal server)" asides
# to the two "Php script to ... the value:" lines, merge the broken-up
# hyperlink text run back into one run, and let the _GoBack bookmark follow
# the latest edit (Word moves it from the old edit location to the new one).

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 0) The _GoBack bookmark currently sits near the end of the document
#    (around the "Screenshot of log ..." paragraph). It is about to be
#    re-created at the new edit location below, and bookmark names must
#    stay unique, so drop the old one first.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1) "Php script to send the value:" -> split "send the value:" into
#    three runs, inserting the parenthetical note before the colon.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("send the value:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target1 = $d.Range($find1.Start, $find1.End)
$target1.InsertXML($pkgHeader + `
    '<w:body><w:p>' + `
    '<w:r><w:t>send the value</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> (note that it is on a custom personal server)</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    '</w:p></w:body>' + $pkgFooter)

# ---------------------------------------------------------------------
# 2) Collapse the three-run hyperlink display text
#    ("http://memes.rel" + "h" + "axmodpack.com/...") back into a single
#    run, leaving the <w:hyperlink> wrapper untouched. A direct
#    Range.Text assignment merges same-format runs in place, but is a
#    no-op when the replacement text is identical to the original, so
#    append a throwaway marker first and then remove it.
# ---------------------------------------------------------------------
$linkUrl = "http://memes.relhaxmodpack.com/SchoolProjects/ELEC3800_WebApp/set_scrollbar_value.php?wrVal=3"
$find2 = $d.Content
$find2.Find.Execute($linkUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find2.Text = $linkUrl + "ZZZ_TMP_MARKER_ZZZ"

$find2b = $d.Content
$find2b.Find.Execute("ZZZ_TMP_MARKER_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find2b.Text = ""

# ---------------------------------------------------------------------
# 3) "Php script to read the value:" -> split into "Php script to read
#    the value" + " " + (new _GoBack bookmark) + "(note that it is on a
#    custom personal server):" - this is where Word's _GoBack bookmark
#    (last edited location) now points.
# ---------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Execute("Php script to read the value:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target3 = $d.Range($find3.Start, $find3.End)
$target3.InsertXML($pkgHeader + `
    '<w:body><w:p>' + `
    '<w:r><w:t>Php script to read the value</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>(note that it is on a custom personal server):</w:t></w:r>' + `
    '</w:p></w:body>' + $pkgFooter)
